$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for StudentId column to include GA
$ws.Range("A1").Value = "StudentId/GA"

# Update GA 6.1, GA 6.2, GA 6.3 scores for the first student row from 4 to 2
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 2
